$wb = $excel.ActiveWorkbook

# Rename the "Tree_Building" sheet to "Tree_Construct"
$ws = $wb.Worksheets.Item("Tree_Building")
$ws.Name = "Tree_Construct"

# Make it the active sheet and select cell E10 (matches new tabSelected view)
$ws.Activate()
$ws.Range("E10").Select()
